# const_infl.xlsx: better dealing with input xlsx -- rename the
# mis-keyed "MEMOXYCHO" species (row 29) to "CH3CHO", and add a sense
# check column (O) that compares the "before"/"after" midday influx
# columns (E and F) for every species row, highlighting rows whose
# diurnal influx profile isn't flat across that boundary.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("const_infl")

# --- touch the font on the rows around the renamed species (A28 HCHO,
#     A30 BENZAL) the way Excel does when you click into neighbouring
#     cells while editing -- this only re-stamps their (unchanged)
#     Calibri font, so visually nothing changes ---
$ws.Range("A28").Font.Name = "Calibri"

# --- rename row 29's species from "MEMOXYCHO" to "CH3CHO" ---
$ws.Range("A29").Style = "Normal"
$ws.Range("A29").Font.Name = "Calibri"
$ws.Range("A29").Font.Size = 12
$ws.Range("A29").Font.Color = 0x1F1F1F
$ws.Range("A29").Value = "CH3CHO"

$ws.Range("A30").Font.Name = "Calibri"

# --- new column O: flag any mismatch between the E (first midday) and
#     F (second midday) influx columns; row 1 (the time-of-day header
#     row) is checked the other way round (E1-F1) ---
$ws.Range("O1").Formula = "=E1-F1"
$ws.Range("O2:O65").Formula = "=F2-E2"
$ws.Range("O66:O76").Formula = "=F66-E66"

# Column P gets a blank, formatted placeholder cell on the header row
# (mirrors the new column O's number format) for symmetry with the
# rest of the header row.
$ws.Range("P1").NumberFormat = "0.0000000"

# --- reposition the view roughly where the author left it ---
[void]$ws.Range("C24").Select()
